# Append 11 new "photo insert SQL" rows (rows 23-33) to the active sheet
# (the "20201026" tab), matching the pattern already used by rows 2-22:
#   - Column A: numeric id
#   - Column B: shared restaurant/shop-id string
#   - Column C: CONCAT(...) formula building the INSERT INTO photos SQL statement
# Odd logical rows (counting from row 2) keep the default style; the others
# reuse the existing "black font" style already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{Row=23; A=11; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$false},
    @{Row=24; A=12; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$true},
    @{Row=25; A=13; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$false},
    @{Row=26; A=14; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$true},
    @{Row=27; A=15; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$false},
    @{Row=28; A=16; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$true},
    @{Row=29; A=17; B="da053f6a-ffb0-11ea-ba65-065a10bcba76"; Styled=$false},
    @{Row=30; A=26; B="da04f9c7-ffb0-11ea-ba65-065a10bcba76"; Styled=$true},
    @{Row=31; A=27; B="da04f9c7-ffb0-11ea-ba65-065a10bcba76"; Styled=$false},
    @{Row=32; A=28; B="da04f9c7-ffb0-11ea-ba65-065a10bcba76"; Styled=$true},
    @{Row=33; A=29; B="da04f9c7-ffb0-11ea-ba65-065a10bcba76"; Styled=$false}
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    $ws.Range("A$r").Value = $rd.A
    $ws.Range("B$r").Value = $rd.B

    $formula = "=_xlfn.CONCAT(""INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('"", B$r, ""'), LPAD("", A$r, "", 7, '0'), 'dish'"", "");"")"
    $ws.Range("C$r").Formula = $formula

    if ($rd.Styled) {
        # Reuse the existing black-font style already used on the
        # alternating rows above (A2, A4, A6, ... A22) instead of the
        # default theme-color font.
        $ws.Range("A$r").Font.Color = 0
    }
}

# Matches the new selection left behind by the edit (B8 instead of B11).
$null = $ws.Range("B8").Select()
